# The "Hofstaat" concordance value "K6" is being renamed to "K6 (Ks.)"
# throughout the table (column D, "Hofstaat"). Scan the used range and
# update every cell whose value is exactly "K6", leaving everything else
# (including the two rows already set to "Dummy Hofstaat") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $used.Row + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $used.Column + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq "K6") {
            $cell.Value = "K6 (Ks.)"
        }
    }
}
